$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# --- New column G: "PRESUPUESTO" budget column, mirroring column F ---

# Header G1 - copy the header formatting (bold, border, centered) from F1,
# then set its text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows G2:G28 - copy the number formatting from F2:F28, then zero them.
$ws.Range("F2:F28").Copy()
$ws.Range("G2:G28").PasteSpecial(-4122)
$ws.Range("G2:G28").Value = 0

# Totals row G29 - copy the totals-row formatting from F29, then zero it.
$ws.Range("F29").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 0

# Column width for the new column G
$ws.Columns.Item(7).ColumnWidth = 16.17

$excel.CutCopyMode = 0
